# Daily attendance processing - reorder the "Recorded By" (column G) author
# lists so each comma-separated list of recorders is written in reverse
# order. Rows whose G cell holds only a single name are left untouched
# (reversing a single-item list is a no-op).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$changedCount = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $orig = $cell.Text

    if ($orig -and $orig.Contains(",")) {
        $parts = $orig -split ', '
        $n = $parts.Count

        $rev = @()
        for ($i = $n - 1; $i -ge 0; $i--) {
            $rev += $parts[$i]
        }
        $newVal = $rev -join ', '

        if ($newVal -ne $orig) {
            $cell.Value = $newVal
            $changedCount++
        }
    }
}

Write-Output ("Reordered Recorded By values on " + $changedCount + " rows")
